$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.856.27"
$ws.Range("E2").Value = "  -2.58%  "
$ws.Range("D3").Value = "2.657.77"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'523.38"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "'144.08"
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  -1.53%  "
$ws.Range("E9").Value = "  +6.77%  "
$ws.Range("D10").Value = "'0.103"
$ws.Range("E10").Value = "  -3.83%  "
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("D13").Value = "3.124.30"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").Value = "58.828.68"
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("D17").Value = "2.664.36"
$ws.Range("E17").Value = "  -6.37%  "
$ws.Range("D18").Value = "'339.13"
$ws.Range("E18").Value = "  -3.68%  "
$ws.Range("E19").Value = "  -3.63%  "
$ws.Range("E20").Value = "  -2.11%  "
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("E27").Value = "  -2.55%  "
$ws.Range("D28").Value = "'7.15"
$ws.Range("E28").Value = "  -2.99%  "
$ws.Range("E29").Value = "  -3.07%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("D32").Value = "'18.88"
$ws.Range("E32").Value = "  -1.61%  "
$ws.Range("D33").Value = "'150.67"
$ws.Range("E33").Value = "  +2.08%  "
$ws.Range("E34").Value = "  -3.74%  "
$ws.Range("B35").Value = "SuiNetwork"
$ws.Range("C35").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D35").Value = "'0.924"
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.19"
$ws.Range("E36").Value = "  -5.14%  "
$ws.Range("D37").Value = "'0.869"
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("D38").Value = "'36.89"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("E39").Value = "  -4.95%  "
$ws.Range("E40").Value = "  -3.38%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").Value = "'275.23"
$ws.Range("E43").Value = "  -4.48%  "
$ws.Range("D44").Value = "'19.72"
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("E45").Value = "  -2.64%  "
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("D47").Value = "2.059.58"
$ws.Range("E47").Value = "  -3.91%  "
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("D49").Value = "'4.72"
$ws.Range("E49").Value = "  -2.63%  "
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("D51").Value = "'18.80"
$ws.Range("E51").Value = "  -2.98%  "
